$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Argon Gas / Base Fit): chi2, Parameters(k), BIC updated
$ws.Range("D2").Value = 3.60949
$ws.Range("E2").Value = 21
$ws.Range("G2").Value = 2555.659236299120

# Row 4 (Argon Gas / 6 Control Points): chi2, Parameters(k), BIC updated
$ws.Range("D4").Value = 4.96996
$ws.Range("E4").Value = 21
$ws.Range("G4").Value = 3467.174136299120

# Row 23 (Magnetic SiOx / Base Fit): chi2, Parameters(k), BIC updated
$ws.Range("D23").Value = 25.9355
$ws.Range("E23").Value = 14
$ws.Range("G23").Value = 9393.746597875010
